# Final commit and update
# - widen column D for the longer product_no text now being stored
# - append new manufacturing data row for g8_testing_T_L8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column D to fit the new, longer product_no values
$ws.Columns.Item(4).ColumnWidth = 17.233072916666668

# Append the new manufacturing data row (row 16)
$ws.Range("A16").Value = 39
$ws.Range("B16").Value = "g8_testing_T_L8"
$ws.Range("C16").Value = "Bottom    "
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = "SMT"
$ws.Range("F16").Value = "SMT_Line_8                    "
$ws.Range("G16").Value = 100
$ws.Range("H16").Value = 43822.0416666667
$ws.Range("H16").NumberFormat = $ws.Range("H2").NumberFormat
$ws.Range("I16").Value = "admin"
